$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Odds updates scattered across existing rows (2,3,4,5,7,8,9) ---
$ws.Cells.Item(2, 13).Value = 1.08  # M2
$ws.Cells.Item(2, 15).Value = 1.44  # O2
$ws.Cells.Item(2, 16).Value = 2.63  # P2
$ws.Cells.Item(2, 22).Value = 1.67  # V2
$ws.Cells.Item(3, 13).Value = 1.03  # M3
$ws.Cells.Item(3, 14).Value = 17  # N3
$ws.Cells.Item(4, 30).Value = 8.5  # AD4
$ws.Cells.Item(5, 7).Value = 2.88  # G5
$ws.Cells.Item(5, 9).Value = 2.55  # I5
$ws.Cells.Item(5, 10).Value = 3.75  # J5
$ws.Cells.Item(5, 12).Value = 3.4  # L5
$ws.Cells.Item(5, 13).Value = 1.11  # M5
$ws.Cells.Item(5, 14).Value = 6.5  # N5
$ws.Cells.Item(5, 23).Value = 7  # W5
$ws.Cells.Item(5, 26).Value = 29  # Z5
$ws.Cells.Item(5, 34).Value = 6.5  # AH5
$ws.Cells.Item(5, 35).Value = 11  # AI5
$ws.Cells.Item(5, 36).Value = 11  # AJ5
$ws.Cells.Item(5, 37).Value = 26  # AK5
$ws.Cells.Item(5, 38).Value = 26  # AL5
$ws.Cells.Item(5, 40).Value = 4.75  # AN5
$ws.Cells.Item(5, 48).Value = 81  # AV5
$ws.Cells.Item(5, 49).Value = 4.5  # AW5
$ws.Cells.Item(5, 53).Value = 101  # BA5
$ws.Cells.Item(7, 13).Value = 1.06  # M7
$ws.Cells.Item(7, 14).Value = 10  # N7
$ws.Cells.Item(8, 7).Value = 1.48  # G8
$ws.Cells.Item(8, 8).Value = 3.8  # H8
$ws.Cells.Item(8, 9).Value = 8  # I8
$ws.Cells.Item(8, 11).Value = 2.2  # K8
$ws.Cells.Item(8, 12).Value = 7.5  # L8
$ws.Cells.Item(8, 21).Value = 2.25  # U8
$ws.Cells.Item(8, 22).Value = 1.57  # V8
$ws.Cells.Item(8, 24).Value = 6  # X8
$ws.Cells.Item(8, 26).Value = 9.5  # Z8
$ws.Cells.Item(8, 35).Value = 41  # AI8
$ws.Cells.Item(8, 36).Value = 23  # AJ8
$ws.Cells.Item(8, 38).Value = 67  # AL8
$ws.Cells.Item(8, 39).Value = 67  # AM8
$ws.Cells.Item(8, 41).Value = 7.5  # AO8
$ws.Cells.Item(8, 43).Value = 23  # AQ8
$ws.Cells.Item(8, 49).Value = 8  # AW8
$ws.Cells.Item(9, 7).Value = 1.75  # G9
$ws.Cells.Item(13, 1).Value = "IyolHyN9"  # A13
$ws.Cells.Item(13, 3).Value = "20:00"  # C13
$ws.Cells.Item(13, 5).Value = "Los Angeles Galaxy"  # E13
$ws.Cells.Item(13, 6).Value = "Minnesota United"  # F13
$ws.Cells.Item(13, 7).Value = 1.73  # G13
$ws.Cells.Item(13, 8).Value = 3.8  # H13
$ws.Cells.Item(13, 9).Value = 4.5  # I13
$ws.Cells.Item(13, 10).Value = 2.25  # J13
$ws.Cells.Item(13, 12).Value = 4.33  # L13
$ws.Cells.Item(13, 14).Value = 17  # N13
$ws.Cells.Item(13, 15).Value = 1.14  # O13
$ws.Cells.Item(13, 16).Value = 5.5  # P13
$ws.Cells.Item(13, 17).Value = 1.53  # Q13
$ws.Cells.Item(13, 18).Value = 2.4  # R13
$ws.Cells.Item(13, 19).Value = 1.25  # S13
$ws.Cells.Item(13, 20).Value = 3.75  # T13
$ws.Cells.Item(13, 21).Value = 1.53  # U13
$ws.Cells.Item(13, 22).Value = 2.38  # V13
$ws.Cells.Item(13, 23).Value = 11  # W13
$ws.Cells.Item(13, 24).Value = 11  # X13
$ws.Cells.Item(13, 26).Value = 15  # Z13
$ws.Cells.Item(13, 28).Value = 19  # AB13
$ws.Cells.Item(13, 29).Value = 17  # AC13
$ws.Cells.Item(13, 31).Value = 12  # AE13
$ws.Cells.Item(13, 32).Value = 34  # AF13
$ws.Cells.Item(13, 33).Value = 101  # AG13
$ws.Cells.Item(13, 35).Value = 26  # AI13
$ws.Cells.Item(13, 36).Value = 15  # AJ13
$ws.Cells.Item(13, 38).Value = 29  # AL13
$ws.Cells.Item(13, 39).Value = 29  # AM13
$ws.Cells.Item(13, 40).Value = 4  # AN13
$ws.Cells.Item(13, 41).Value = 9  # AO13
$ws.Cells.Item(13, 43).Value = 26  # AQ13
$ws.Cells.Item(13, 47).Value = 7.5  # AU13
$ws.Cells.Item(13, 48).Value = 41  # AV13
$ws.Cells.Item(13, 49).Value = 6.5  # AW13
$ws.Cells.Item(13, 50).Value = 21  # AX13
$ws.Cells.Item(13, 51).Value = 23  # AY13
$ws.Cells.Item(13, 52).Value = 67  # AZ13
$ws.Cells.Item(13, 53).Value = 67  # BA13
$ws.Cells.Item(13, 54).Value = 126  # BB13
$ws.Cells.Item(13, 55).Value = 351  # BC13

# --- New match row 14: Dep. Tachira vs Carabobo (VENEZUELA - LIGA FUTVE) ---
$ws.Cells.Item(14, 1).Value = "l0F5POqe"  # A14
$ws.Cells.Item(14, 2).Value = "24/11/2024"  # B14
$ws.Cells.Item(14, 3).Value = "20:00"  # C14
$ws.Cells.Item(14, 4).Value = "VENEZUELA - LIGA FUTVE"  # D14
$ws.Cells.Item(14, 5).Value = "Dep. Tachira"  # E14
$ws.Cells.Item(14, 6).Value = "Carabobo"  # F14
$ws.Cells.Item(14, 7).Value = 1.65  # G14
$ws.Cells.Item(14, 8).Value = 3.2  # H14
$ws.Cells.Item(14, 9).Value = 5.8  # I14
$ws.Cells.Item(14, 10).Value = 2.27  # J14
$ws.Cells.Item(14, 11).Value = 1.95  # K14
$ws.Cells.Item(14, 12).Value = 6  # L14
$ws.Cells.Item(14, 13).Value = 1.07  # M14
$ws.Cells.Item(14, 14).Value = 5.45  # N14
$ws.Cells.Item(14, 15).Value = 1.5  # O14
$ws.Cells.Item(14, 16).Value = 2.27  # P14
$ws.Cells.Item(14, 17).Value = 2.4  # Q14
$ws.Cells.Item(14, 18).Value = 1.44  # R14
$ws.Cells.Item(14, 19).Value = 1.52  # S14
$ws.Cells.Item(14, 20).Value = 2.22  # T14
$ws.Cells.Item(14, 21).Value = 2.27  # U14
$ws.Cells.Item(14, 22).Value = 1.5  # V14
$ws.Cells.Item(14, 23).Value = 4.65  # W14
$ws.Cells.Item(14, 24).Value = 6.2  # X14
$ws.Cells.Item(14, 25).Value = 9  # Y14
$ws.Cells.Item(14, 26).Value = 12  # Z14
$ws.Cells.Item(14, 27).Value = 17.5  # AA14
$ws.Cells.Item(14, 28).Value = 45  # AB14
$ws.Cells.Item(14, 29).Value = 6.3  # AC14
$ws.Cells.Item(14, 30).Value = 6.7  # AD14
$ws.Cells.Item(14, 31).Value = 24  # AE14
$ws.Cells.Item(14, 32).Value = 175  # AF14
$ws.Cells.Item(14, 33).Value = 101  # AG14
$ws.Cells.Item(14, 34).Value = 11  # AH14
$ws.Cells.Item(14, 35).Value = 32  # AI14
$ws.Cells.Item(14, 36).Value = 20  # AJ14
$ws.Cells.Item(14, 37).Value = 150  # AK14
$ws.Cells.Item(14, 38).Value = 90  # AL14
$ws.Cells.Item(14, 39).Value = 100  # AM14
$ws.Cells.Item(14, 40).Value = 3.2  # AN14
$ws.Cells.Item(14, 41).Value = 8.25  # AO14
$ws.Cells.Item(14, 42).Value = 23  # AP14
$ws.Cells.Item(14, 43).Value = 30  # AQ14
$ws.Cells.Item(14, 44).Value = 90  # AR14
$ws.Cells.Item(14, 45).Value = 400  # AS14
$ws.Cells.Item(14, 46).Value = 2.2  # AT14
$ws.Cells.Item(14, 47).Value = 8.75  # AU14
$ws.Cells.Item(14, 48).Value = 110  # AV14
$ws.Cells.Item(14, 49).Value = 6.9  # AW14
$ws.Cells.Item(14, 50).Value = 37  # AX14
$ws.Cells.Item(14, 51).Value = 45  # AY14
$ws.Cells.Item(14, 52).Value = 300  # AZ14
$ws.Cells.Item(14, 53).Value = 350  # BA14
